$d = $word.ActiveDocument

# The entire document uses a uniform font size of 9pt (sz/szCs = 18 half-points).
# Bump it to 10pt (sz/szCs = 20 half-points) across every paragraph and run,
# including the paragraph-mark run properties. Font.Size maps to w:sz while
# Font.SizeBi maps to w:szCs (complex-script size), so both need to be set.
$d.Content.Font.Size = 10
$d.Content.Font.SizeBi = 10

foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Size = 10
    $p.Range.Font.SizeBi = 10
}
